$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C2").Value = 3300.0
$ws.Range("C4").Value = 106.0
$ws.Range("C6").Value = 68885.4928957218
$ws.Range("C7").Value = 64645.49289572181
$ws.Range("C8").Value = 66818.92810885014
$ws.Range("C9").Value = 9540.0
$ws.Range("C10").Value = 17638.63418025324
$ws.Range("C11").Value = 23212.319159666004
$ws.Range("C13").Value = 45673.17373605581
$ws.Range("C14").Value = 41433.17373605581
$ws.Range("C15").Value = 36133.17373605582
$ws.Range("C16").Value = 35559.202614714595
$ws.Range("C17").Value = 34556.33009055581
$ws.Range("C19").Value = 344.4274758412121
$ws.Range("C21").Value = 675535.91890583
$ws.Range("C22").Value = 633955.72290583
$ws.Range("C23").Value = 655269.841338655
$ws.Range("C24").Value = 93555.44099999996
$ws.Range("C25").Value = 172975.91188378038
$ws.Range("C27").Value = 447900.82921869156
$ws.Range("C28").Value = 406320.63321869157
$ws.Range("C29").Value = 354345.3882186917
$ws.Range("C30").Value = 348716.6543215908
$ws.Range("C31").Value = 338881.8344825491
$ws.Range("C33").Value = 3377.6797059582213

$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C6").Value = 9743.0
$ws.Range("D6").Value = 9.625879043600607
$ws.Range("C7").Value = 7210.0
$ws.Range("D7").Value = -18.874824191279856
$ws.Range("C8").Value = 6703.0
$ws.Range("D8").Value = -24.57946554149083
$ws.Range("C9").Value = 7121.0
$ws.Range("D9").Value = -19.876230661040758
$ws.Range("C12").Value = 7811.333333333332
$ws.Range("D12").Value = -12.108766994842924

$ws = $wb.Worksheets.Item("WING")
$ws.Range("C3").Value = 1.0
$ws.Range("C7").Value = 0.0
$ws.Range("D7").Value = -99.99999999999999
$ws.Range("C8").Value = 6387.0
$ws.Range("D8").Value = 7.797468354430412
$ws.Range("C10").Value = 7361.0
$ws.Range("D10").Value = 24.23628691983126
$ws.Range("C14").Value = 5765.285714285714
$ws.Range("D14").Value = -2.695599758890876

$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C3").Value = 1.0
$ws.Range("C8").Value = 276.0
$ws.Range("D8").Value = -68.94514767932489
$ws.Range("C9").Value = 556.0
$ws.Range("D9").Value = -37.44022503516173
$ws.Range("C10").Value = 485.0
$ws.Range("D10").Value = -45.428973277074526

$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C3").Value = 1.0

$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C3").Value = 893.3333333333333
$ws.Range("D3").Value = -66.49476480700108
$ws.Range("C9").Value = 225.0
$ws.Range("D9").Value = -49.36708860759492
$ws.Range("C10").Value = 289.0
$ws.Range("D10").Value = -34.96483825597748
$ws.Range("C11").Value = 826.0
$ws.Range("D11").Value = 85.87904360056264
$ws.Range("C12").Value = 446.66666666666663
$ws.Range("C16").Value = 225.0
$ws.Range("D16").Value = -49.36708860759492
$ws.Range("C17").Value = 289.0
$ws.Range("D17").Value = -34.96483825597748
$ws.Range("C18").Value = 826.0
$ws.Range("D18").Value = 85.87904360056264
$ws.Range("C19").Value = 446.66666666666663

$ws = $wb.Worksheets.Item("POWER PLANT")
$ws.Range("C3").Value = 9140.767626114684
$ws.Range("C9").Value = 4570.383813057342
$ws.Range("C13").Value = 4570.383813057342

$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C5").Value = 1918.0
$ws.Range("D5").Value = -19.07172995780588
$ws.Range("C6").Value = 2755.0
$ws.Range("D6").Value = 16.244725738396667
$ws.Range("C7").Value = 3183.0
$ws.Range("D7").Value = 34.30379746835448
$ws.Range("C8").Value = 2791.0
$ws.Range("D8").Value = 17.76371308016882
$ws.Range("C9").Value = 2661.75
$ws.Range("D9").Value = 12.310126582278501
